$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("學生名單")

# Insert 8 new rows below the existing data (rows 2-7), copying row 7's
# formatting so the new rows match the existing table style.
for ($i = 0; $i -lt 8; $i++) {
    $row = 8 + $i
    $ws.Rows(7).Copy()
    $ws.Rows($row).Insert()
}

$names = @("呂彥臻1", "王大同1", "楊小明1", "劉大象1", "陳中一1", "孫二1", "呂彥臻2", "王大同2")
$classes = @("一", "二", "三", "四", "五", "六", "一", "二")
$grades = @("忠", "孝", "仁", "愛", "信", "義", "忠", "孝")

for ($i = 0; $i -lt 8; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $classes[$i]
    $ws.Cells.Item($row, 3).Value = $grades[$i]
}

$ws.Activate()
$ws.Range("A15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4 | Out-Null
